$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.217.92'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.521.13'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.62'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.36'
$ws.Range('E6').Value = '  -3.82%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('D9').Value = '2.527.32'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.44'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.356'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '2.967.62'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.30'
$ws.Range('E15').Value = '  -2.51%  '
$ws.Range('D16').Value = '59.150.38'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000141'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '2.520.13'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.05'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.31'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.52'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.83'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.40'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.427'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.82'
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.78'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('D30').Value = '0.0₃0774'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E31').Value = '  -2.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.95'
$ws.Range('E32').Value = '  +4.91%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.47'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('E35').Value = '  -8.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.52'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.26'
$ws.Range('E37').Value = '  -3.49%  '
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.92'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.68'
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'
$ws.Range('E41').Value = '  -2.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.26'
$ws.Range('E42').Value = '  -7.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '279.85'
$ws.Range('E43').Value = '  -6.54%  '
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.85'
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0934'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.52'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.53'
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0513'
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('E51').Value = '  -2.22%  '
